$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "326.51", "1.120").
# Excel auto-converts such literals to real numbers on assignment, which
# rounds/truncates them (floating point) and loses trailing zeros. Forcing
# the whole D2:D51 range to Text format before writing keeps them as exact
# strings; resetting the style back to Normal afterwards avoids leaving any
# stray number-format styling behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.191.82"
$ws.Range("D3").Value = "1.904.02"
$ws.Range("D5").Value = "326.51"
$ws.Range("D7").Value = "0.5158"
$ws.Range("D8").Value = "0.4009"
$ws.Range("D9").Value = "0.08465"
$ws.Range("D10").Value = "42.63"
$ws.Range("D11").Value = "1.120"
$ws.Range("D12").Value = "23.26"
$ws.Range("D13").Value = "6.437"
$ws.Range("D14").Value = "1.904.81"
$ws.Range("D15").Value = "7.353"
$ws.Range("D17").Value = "94.86"
$ws.Range("D19").Value = "0.06648"
$ws.Range("D20").Value = "18.36"
$ws.Range("D22").Value = "5.997"
$ws.Range("D23").Value = "30.199.81"
$ws.Range("D24").Value = "11.26"
$ws.Range("D26").Value = "2.127.56"
$ws.Range("D27").Value = "21.59"
$ws.Range("D28").Value = "162.15"
$ws.Range("D29").Value = "2.386"
$ws.Range("D30").Value = "129.32"
$ws.Range("D31").Value = "1.097"
$ws.Range("D33").Value = "6.079"
$ws.Range("D34").Value = "3.679"
$ws.Range("D35").Value = "0.02490"
$ws.Range("D36").Value = "0.06566"
$ws.Range("D37").Value = "0.2201"
$ws.Range("D38").Value = "5.196"
$ws.Range("D39").Value = "1.227"
$ws.Range("D40").Value = "11.95"
$ws.Range("D41").Value = "8.795"
$ws.Range("D42").Value = "0.6513"
$ws.Range("D43").Value = "1.233"
$ws.Range("D44").Value = "0.6124"
$ws.Range("D45").Value = "13.27"
$ws.Range("D48").Value = "1.244"
$ws.Range("D49").Value = "125.02"
$ws.Range("D50").Value = "1.163"
$ws.Range("D51").Value = "79.15"

$dRange.Style = "Normal"

# Column E holds percentage strings with padding spaces, e.g. "  +3.36%  ".
# These are not parsed as numbers by Excel (leading/trailing spaces), so
# they can be assigned directly.
$ws.Range("E2").Value = "  +3.36%  "
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +3.43%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("E12").Value = "  +13.26%  "
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  +3.43%  "
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("E31").Value = "  +3.32%  "
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  +5.75%  "
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  +1.82%  "
